# Update the LR-pairs (Adm-Calcrl) sheet with the new TPM-normalized values.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity (avg), J Ligand derived specificity (total),
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity (avg), P Receptor derived specificity (total),
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ G=15.89577633333333; H=47.687329;  I=0.286059172443548;  J=0.2860591724435479;
           M=33.61498366666667; N=100.844951; O=0.4159547139783538; P=0.4159547139783538;
           Q=534.3362618139865; R=4809.02635632588; S=0.1189876612546406; T=0.1189876612546406 }
    3 = @{ G=15.89577633333333; H=47.687329;  I=0.286059172443548;  J=0.2860591724435479;
           O=0.4210695178651698; P=0.4210695178651698;
           Q=540.9067491698172; R=4868.160742528355; S=0.1204507978217142; T=0.1204507978217142 }
    4 = @{ G=15.89577633333333; H=47.687329;  I=0.286059172443548;  J=0.2860591724435479;
           O=0.1629757681564764; P=0.1629757681564764;
           Q=209.3589994210918; R=1884.230994789826; S=0.04662071336719319; T=0.04662071336719319 }
    5 = @{ I=0.6735478078679881; J=0.673547807867988;
           M=33.61498366666667; N=100.844951; O=0.4159547139783538; P=0.4159547139783538;
           Q=1258.134863269278; R=11323.2137694235; S=0.2801653857724761; T=0.2801653857724761 }
    6 = @{ I=0.6735478078679881; J=0.673547807867988;
           O=0.4210695178651698; P=0.4210695178651698;
           S=0.2836104507181157; T=0.2836104507181157 }
    7 = @{ I=0.6735478078679881; J=0.673547807867988;
           O=0.1629757681564764; P=0.1629757681564764;
           S=0.1097719713773962; T=0.1097719713773961 }
    8 = @{ I=0.04039301968846393; J=0.04039301968846393;
           M=33.61498366666667; N=100.844951; O=0.4159547139783538; P=0.4159547139783538;
           Q=75.45101581377168; R=679.059142323945; S=0.01680166695123703; T=0.01680166695123703 }
    9 = @{ I=0.04039301968846393; J=0.04039301968846393;
           O=0.4210695178651698; P=0.4210695178651698;
           S=0.01700826932533982; T=0.01700826932533982 }
    10 = @{ I=0.04039301968846393; J=0.04039301968846393;
           O=0.1629757681564764; P=0.1629757681564764;
           S=0.006583083411887086; T=0.006583083411887085 }
}

foreach ($rowNum in $data.Keys) {
    $cols = $data[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
